$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "27.279.78"
Set-TextCell "E2" "  -2.53%  "

# Row 3
Set-TextCell "D3" "1.706.84"
Set-TextCell "E3" "  -2.17%  "

# Row 4
Set-TextCell "E4" "  -0.02%  "

# Row 5
Set-TextCell "D5" "223.38"
Set-TextCell "E5" "  -2.10%  "

# Row 6
Set-TextCell "D6" "0.5316"
Set-TextCell "E6" "  -2.17%  "

# Row 7
Set-TextCell "D7" "1.002"
Set-TextCell "E7" "  -0.03%  "

# Row 8
Set-TextCell "D8" "0.2657"
Set-TextCell "E8" "  -3.91%  "

# Row 9
Set-TextCell "D9" "0.06584"
Set-TextCell "E9" "  -2.34%  "

# Row 10
Set-TextCell "D10" "20.81"
Set-TextCell "E10" "  -4.10%  "

# Row 11
Set-TextCell "D11" "0.07635"
Set-TextCell "E11" "  -1.91%  "

# Row 12
Set-TextCell "D12" "4.562"
Set-TextCell "E12" "  -3.13%  "

# Row 13
Set-TextCell "D13" "1.710.69"
Set-TextCell "E13" "  -1.89%  "

# Row 14
Set-TextCell "D14" "1.944.71"
Set-TextCell "E14" "  -1.89%  "

# Row 15
Set-TextCell "D15" "0.5720"
Set-TextCell "E15" "  -4.34%  "

# Row 16
Set-TextCell "D16" "0.0₅8163"
Set-TextCell "E16" "  -2.64%  "

# Row 17
Set-TextCell "D17" "67.70"
Set-TextCell "E17" "  -1.83%  "

# Row 18
Set-TextCell "D18" "27.338.90"
Set-TextCell "E18" "  -2.28%  "

# Row 19
Set-TextCell "D19" "215.96"
Set-TextCell "E19" "  -3.81%  "

# Row 20
Set-TextCell "D20" "1.002"
Set-TextCell "E20" "  -0.06%  "

# Row 21
Set-TextCell "D21" "4.663"
Set-TextCell "E21" "  -3.77%  "

# Row 22
Set-TextCell "D22" "10.42"
Set-TextCell "E22" "  -4.46%  "

# Row 23
Set-TextCell "D23" "5.965"
Set-TextCell "E23" "  -4.40%  "

# Row 24
Set-TextCell "D24" "1.003"
Set-TextCell "E24" "  -0.04%  "

# Row 25
Set-TextCell "D25" "1.767"
Set-TextCell "E25" "  +5.70%  "

# Row 26
Set-TextCell "D26" "141.56"
Set-TextCell "E26" "  -3.23%  "

# Row 27
Set-TextCell "D27" "0.1214"
Set-TextCell "E27" "  -2.54%  "

# Row 28
Set-TextCell "D28" "7.259"
Set-TextCell "E28" "  -2.66%  "

# Row 29
Set-TextCell "D29" "16.29"
Set-TextCell "E29" "  -5.94%  "

# Row 30
Set-TextCell "D30" "0.05415"
Set-TextCell "E30" "  -4.05%  "

# Row 31
Set-TextCell "D31" "1.291"
Set-TextCell "E31" "  -1.75%  "

# Row 32
Set-TextCell "D32" "3.500"
Set-TextCell "E32" "  -5.58%  "

# Row 33
Set-TextCell "D33" "3.423"
Set-TextCell "E33" "  -2.77%  "

# Row 34
Set-TextCell "E34" "  -2.26%  "

# Row 35
Set-TextCell "D35" "2.874"
Set-TextCell "E35" "  +0.51%  "

# Row 36
Set-TextCell "D36" "0.9476"
Set-TextCell "E36" "  -3.50%  "

# Row 37
Set-TextCell "D37" "2.406"
Set-TextCell "E37" "  -1.77%  "

# Row 38
Set-TextCell "D38" "0.5855"
Set-TextCell "E38" "  -1.92%  "

# Row 39
Set-TextCell "D39" "0.01629"
Set-TextCell "E39" "  -2.14%  "

# Row 40
Set-TextCell "D40" "5.856"
Set-TextCell "E40" "  -1.47%  "

# Row 41
Set-TextCell "D41" "1.043.05"
Set-TextCell "E41" "  -0.65%  "

# Row 42
Set-TextCell "B42" "TrustWalletToken"
Set-TextCell "C42" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D42" "0.8430"
Set-TextCell "E42" "  -0.86%  "

# Row 43
Set-TextCell "B43" "PaxDollar"
Set-TextCell "C43" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell "D43" "1.002"
Set-TextCell "E43" "  -0.07%  "

# Row 44
Set-TextCell "D44" "100.77"
Set-TextCell "E44" "  -1.23%  "

# Row 45
Set-TextCell "D45" "1.850.52"
Set-TextCell "E45" "  -1.97%  "

# Row 46
Set-TextCell "D46" "0.0₈113"
Set-TextCell "E46" "  -1.82%  "

# Row 47
Set-TextCell "D47" "57.90"
Set-TextCell "E47" "  -3.49%  "

# Row 48
Set-TextCell "D48" "0.4500"
Set-TextCell "E48" "  +1.51%  "

# Row 49
Set-TextCell "D49" "1.002"
Set-TextCell "E49" "  -0.08%  "

# Row 50
Set-TextCell "D50" "8.070"
Set-TextCell "E50" "  -2.21%  "

# Row 51
Set-TextCell "D51" "0.05242"
Set-TextCell "E51" "  -1.55%  "
